# Update the "Abnormal Events" sheet with refreshed vital-sign anomaly data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Abnormal Events")

$data = @(
    @("2026-01-13 21:35:05", "2026-01-13 21:35:05", "ECG",         "-0.9 - -0.9", "AMBER"),
    @("2026-01-13 21:35:12", "2026-01-13 21:35:12", "ECG",         "-0.8 - -0.8", "AMBER"),
    @("2026-01-13 21:35:05", "2026-01-13 21:35:16", "Temperature", "35.1 - 35.9", "AMBER"),
    @("2026-01-13 21:35:19", "2026-01-13 21:35:19", "ECG",         "0.9 - 0.9",   "RED"),
    @("2026-01-13 21:35:20", "2026-01-13 21:35:21", "ECG",         "-0.9 - -0.8", "AMBER")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
